$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Eddie"
$ws.Range("B2").Value = "Espinoza"
$ws.Range("C2").Value = 75
